$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J1/K1 used to hold shared strings "r"/"s" - replace with numeric 0.5 each
$ws.Range("J1").Value = 0.5
$ws.Range("K1").Value = 0.5

# J2:J51 change from 0.6 to 0.5 (K2:K51 already 0.5, left untouched)
$ws.Range("J2:J51").Value = 0.5

# Update the visible selection to match the retrained-model view: K1:K51, active cell K1
$ws.Range("K1:K51").Select()
